$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2/E2 values (shared string content changes: testT -> test11, recT -> sender11)
$ws.Range("D2").Value = "test11"
$ws.Range("E2").Value = "sender11"

# Update D3/E3 values (testR -> test10, recR -> mem)
$ws.Range("D3").Value = "test10"
$ws.Range("E3").Value = "mem"

# Update A2 and A3 values 2 -> 3
$ws.Range("A2").Value = 3
$ws.Range("A3").Value = 3

# Add new row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "Dropoff"
$ws.Range("D4").Value = "test12"
$ws.Range("E4").Value = "rec12"

# Update selection to E2
$ws.Range("E2").Select()
